$d = $word.ActiveDocument

# 1) Merge the two runs that make up the orderDoctor field placeholder.
#    The original text is split across two runs with a "_GoBack" bookmark
#    sitting between them; replacing the full phrase in one Find/Replace
#    call merges it back into a single run and drops the bookmark.
$d.Content.Find.Execute(
    "[report.orderDoctor; block=tbs:row+tbs:row+tbs:row; sub1=transactions; p1] [report.orderDoctorNameTH;]",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[report.orderDoctor; block=tbs:row+tbs:row+tbs:row; sub1=transactions; p1] [report.orderDoctorNameTH;]",
    2) | Out-Null

# 2) Update the insurance field placeholder text (stays within the same run).
$d.Content.Find.Execute(
    "[report_sub1.insurance.condition.insuranceName;ifempty=",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "[report_sub1.insurance.PatientsInsurances;ope=formatinsurance;ifempty=",
    2) | Out-Null

# 3) Re-insert the "_GoBack" bookmark around the whole insurance placeholder
#    field (from the start of its first run "[report_sub1..." to the end of
#    its closing "]" run). Locate the start and end independently so this
#    does not depend on the embedded (language-specific) default-value text.
$startRng = $d.Content
$startRng.Find.Execute(
    "[report_sub1.insurance.PatientsInsurances;ope=formatinsurance;ifempty=",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$fieldStart = $startRng.Start

$endRng = $d.Content
$endRng.Start = $fieldStart
$endRng.Find.Execute("]", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null
$fieldEnd = $endRng.End

$rng = $d.Range($fieldStart, $fieldEnd)
$d.Bookmarks.Add("_GoBack", $rng)
